$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 76923280
$ws.Range("J55").Value = 165.2
$ws.Range("L55").Value = 165.2
$ws.Range("N55").Value = -593.2
$ws.Range("H92").Value = 352.8
$ws.Range("I92").Value = 116.333336
$ws.Range("K92").Value = 116.333336
$ws.Range("M92").Value = 1131.666664
$ws.Range("H107").Value = 771.44446
$ws.Range("I107").Value = 771.44446
$ws.Range("K107").Value = 771.44446
$ws.Range("M107").Value = 1148.55554
$ws.Range("H134").Value = 38995.668
$ws.Range("J134").Value = 38995.668
$ws.Range("L134").Value = 38995.668
$ws.Range("N134").Value = -49135.668

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1654.5883
$ws.Range("J2").Value = 2131.5
$ws.Range("L2").Value = 2131.5
$ws.Range("N2").Value = -2357.5
$ws.Range("H32").Value = 8026.742
$ws.Range("I32").Value = 4850.625
$ws.Range("K32").Value = 4850.625
$ws.Range("M32").Value = -4563.625
$ws.Range("H45").Value = 2688.7273
$ws.Range("I45").Value = 1200
$ws.Range("K45").Value = 1200
$ws.Range("M45").Value = -823
$ws.Range("H74").Value = 2422.8333
$ws.Range("I74").Value = 1280
$ws.Range("K74").Value = 1280
$ws.Range("M74").Value = -406
$ws.Range("H77").Value = 2422.8333
$ws.Range("I77").Value = 1280
$ws.Range("K77").Value = 6400
$ws.Range("M77").Value = -2032
$ws.Range("H88").Value = 8824.177
$ws.Range("I88").Value = 1340
$ws.Range("J88").Value = 11127
$ws.Range("K88").Value = 1340
$ws.Range("L88").Value = 11127
$ws.Range("M88").Value = -934
$ws.Range("N88").Value = -11939
$ws.Range("H91").Value = 8824.177
$ws.Range("I91").Value = 1340
$ws.Range("J91").Value = 11127
$ws.Range("K91").Value = 1340
$ws.Range("L91").Value = 11127
$ws.Range("M91").Value = 64
$ws.Range("N91").Value = -13935
$ws.Range("H116").Value = 1654.5883
$ws.Range("J116").Value = 2131.5
$ws.Range("L116").Value = 2131.5
$ws.Range("N116").Value = -6719.5
$ws.Range("H122").Value = 7016.393
$ws.Range("I122").Value = 7094.077
$ws.Range("J122").Value = 6006.5
$ws.Range("K122").Value = 21282.231
$ws.Range("L122").Value = 18019.5
$ws.Range("M122").Value = -18832.231
$ws.Range("N122").Value = -22919.5
$ws.Range("H132").Value = 2923.7083
$ws.Range("I132").Value = 2755.8096
$ws.Range("K132").Value = 8267.4288
$ws.Range("M132").Value = -5737.4288

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1654.5883
$ws.Range("J3").Value = 2131.5
$ws.Range("L3").Value = 2131.5
$ws.Range("N3").Value = -2359.5
$ws.Range("H17").Value = 6002.5
$ws.Range("J17").Value = 6002.5
$ws.Range("L17").Value = 6002.5
$ws.Range("N17").Value = -6346.5
$ws.Range("H94").Value = 12757.75
$ws.Range("I94").Value = 282.76923
$ws.Range("K94").Value = 282.76923
$ws.Range("M94").Value = 168.23077
$ws.Range("H105").Value = 3906.4
$ws.Range("I105").Value = 2953
$ws.Range("K105").Value = 2953
$ws.Range("M105").Value = -1206
$ws.Range("H133").Value = 75396.336
$ws.Range("J133").Value = 75396.336
$ws.Range("L133").Value = 75396.336
$ws.Range("N133").Value = -85516.336

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 755.64703
$ws.Range("I22").Value = 632.25
$ws.Range("K22").Value = 632.25
$ws.Range("M22").Value = -282.25
$ws.Range("H132").Value = 2273.0322
$ws.Range("I132").Value = 2273.0322
$ws.Range("K132").Value = 6819.096600000001
$ws.Range("M132").Value = -4289.096600000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 2246.5
$ws.Range("I86").Value = 3994
$ws.Range("J86").Value = 499
$ws.Range("K86").Value = 11982
$ws.Range("L86").Value = 1497
$ws.Range("M86").Value = -10796
$ws.Range("N86").Value = -3869
$ws.Range("H89").Value = 2246.5
$ws.Range("I89").Value = 3994
$ws.Range("J89").Value = 499
$ws.Range("K89").Value = 35946
$ws.Range("L89").Value = 4491
$ws.Range("M89").Value = -30018
$ws.Range("N89").Value = -16347
$ws.Range("H131").Value = 14410290
$ws.Range("J131").Value = 13941749
$ws.Range("L131").Value = 41825247
$ws.Range("N131").Value = -41835327

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3409.182
$ws.Range("J80").Value = 3324.5
$ws.Range("L80").Value = 3324.5
$ws.Range("N80").Value = -5320.5
$ws.Range("H83").Value = 3409.182
$ws.Range("J83").Value = 3324.5
$ws.Range("L83").Value = 16622.5
$ws.Range("N83").Value = -26606.5
$ws.Range("H102").Value = 13589.637
$ws.Range("I102").Value = 15776.889
$ws.Range("J102").Value = 3747
$ws.Range("K102").Value = 15776.889
$ws.Range("L102").Value = 3747
$ws.Range("M102").Value = -14154.889
$ws.Range("N102").Value = -6991
$ws.Range("H113").Value = 4111
$ws.Range("I113").Value = 4358
$ws.Range("J113").Value = 3246.5
$ws.Range("K113").Value = 4358
$ws.Range("L113").Value = 3246.5
$ws.Range("M113").Value = -2188
$ws.Range("N113").Value = -7586.5
$ws.Range("H126").Value = 7295.125
$ws.Range("I126").Value = 6302.933
$ws.Range("K126").Value = 18908.799
$ws.Range("M126").Value = -16438.799
$ws.Range("H135").Value = 97998
$ws.Range("J135").Value = 97998
$ws.Range("L135").Value = 97998
$ws.Range("N135").Value = -108138

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1926.6875
$ws.Range("I16").Value = 1626.0714
$ws.Range("J16").Value = 4031
$ws.Range("K16").Value = 1626.0714
$ws.Range("L16").Value = 4031
$ws.Range("M16").Value = -1456.0714
$ws.Range("N16").Value = -4371
$ws.Range("H40").Value = 2779.8
$ws.Range("I40").Value = 2824.5
$ws.Range("K40").Value = 2824.5
$ws.Range("M40").Value = -2688.5
$ws.Range("H46").Value = 2232.7778
$ws.Range("I46").Value = 1388.2
$ws.Range("K46").Value = 1388.2
$ws.Range("M46").Value = -1200.2
$ws.Range("H94").Value = 40165
$ws.Range("J94").Value = 40165
$ws.Range("L94").Value = 40165
$ws.Range("N94").Value = -41517
$ws.Range("H136").Value = 3280.761
$ws.Range("I136").Value = 3196.4
$ws.Range("K136").Value = 9589.200000000001
$ws.Range("M136").Value = -7039.200000000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 27432.666
$ws.Range("J101").Value = 27432.666
$ws.Range("L101").Value = 27432.666
$ws.Range("N101").Value = -33922.666
$ws.Range("H122").Value = 334433
$ws.Range("I122").Value = 334433
$ws.Range("K122").Value = 1003299
$ws.Range("M122").Value = -1000849
$ws.Range("H132").Value = 5709.724
$ws.Range("I132").Value = 7284.15
$ws.Range("K132").Value = 21852.45
$ws.Range("M132").Value = -19322.45
